# Add a new "VLAN_error" column to the global_test sheet.
#
# The source sheet has a header row (row 1) and one data row (row 2) laid
# out from column A to column AD. A new column is inserted right before
# column L ("VLAN#Data" / Trunk group), pushing every following column one
# slot to the right (L:AD -> M:AE). The new column gets the header
# "VLAN_error" and a numeric sample value of 109 for the single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at L, shifting L:AD to M:AE (keeps formatting/styles
# of the column that used to sit at L, matching Excel's normal Insert
# behaviour).
$null = $ws.Columns("L:L").Insert()

# Populate the newly inserted column.
$ws.Range("L1").Value = "VLAN_error"
$ws.Range("L2").Value = 109

# Leave the selection where the author ended up after the edit.
$null = $ws.Range("L3").Select()
